$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert two new columns before the old column B ("Unnamed: 0.18"),
#    shifting all existing columns B..V right to D..X.
# ---------------------------------------------------------------
$ws.Columns("B:C").Insert()

# The Insert() call copies formatting from the column to the left (A),
# which is bold/bordered/centered. Data rows (2..111) in the two new
# columns B:C must NOT have that header styling - only column A and
# row 1 keep it. Clear it from the data rows now; the new B1/C1 header
# cells get the header styling applied further down.
$ws.Range("B2:C111").Font.Bold = $false
$ws.Range("B2:C111").Borders.LineStyle = 0
$ws.Range("B2:C111").HorizontalAlignment = -4131
$ws.Range("B2:C111").VerticalAlignment = -4160

# ---------------------------------------------------------------
# 2. Header row: write the two new "Unnamed: 0.x" labels and style them
#    like the rest of row 1 (bold, thin border, centered/top aligned).
# ---------------------------------------------------------------
$ws.Range("B1").Value = "Unnamed: 0.20"
$ws.Range("C1").Value = "Unnamed: 0.19"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").Borders.LineStyle = 1
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4160

# ---------------------------------------------------------------
# 3. Data: this sheet is a pandas "staircase" export - 22 index columns
#    (A..V), where column j (0 = A) holds the row's running index n for
#    as long as (n - j) < 5 * (22 - j)/5 ... i.e. the last populated
#    column index for row n is (21 - floor(n/5)). Rebuild rows 2..111
#    (n = 0..109) to follow that rule across all 22 index columns, and
#    fill the two value columns W ("   SILVER_FOR", raw) and
#    X ("SILVER_FOR", final) from the known source data.
# ---------------------------------------------------------------
$rawValues = @(27.44466,26.933548,27.173124,26.856646,26.422922,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
$finalValues = @($null,$null,$null,$null,$null,30.93059290717292,30.6353459147918,30.30023174736436,29.53532024840348,29.15916464141611,31.76244298992617,32.23934168175691,32.5524573182023,32.38868114596181,32.10415671664822,28.81685345771996,28.51061683625062,27.67123246313037,27.34084839285276,26.94339908891732,29.09709106441289,28.96379896414152,28.94952882821161,29.52255201468267,28.89384280964123,28.95396845664823,28.94541683747048,28.70188153581586,27.87425082086361,28.09067066164505,31.23943548728738,31.04450482476813,30.54658087863896,30.17317577633884,29.77593984989414,29.80264645196246,29.55922128062565,29.30527315742233,28.8122954391917,28.18952363490905,30.70601929298681,30.43889270062886,30.30290516082914,30.37639649853105,29.82201011509096,28.05217323242647,27.89165416433514,28.205499414136,28.41437395151866,28.78275573469779,31.55722890898164,30.83284674360289,30.74265827309239,30.71705159757317,30.40893174766933,30.00756593430685,29.61517987173261,28.95441448443535,28.93568544726321,29.07687113806344,30.66437626656585,29.71554851349543,29.37691292915611,29.73525220398017,30.20347400339665,30.05011065270486,29.4726120454439,29.81097643387562,29.90499743944582,29.14887020995906,30.23977983929086,30.38220174290848,30.49393445684007,30.98144420479156,30.88058823637243,29.5957914292826,29.50264458311,29.43434864315066,29.30550030120281,28.88129208329809,30.46725951723812,30.30320113605825,30.33314281621216,30.25661876488363,30.87639364789294,28.66571473944788,28.8619683424098,29.76204140378684,29.70551151165625,29.53626048599176,31.81960475475518,31.37961348567967,30.52227511749885,30.01177979559196,29.68068125487758,28.48524969417264,28.19800917589612,28.07364451451201,27.45626358928604,27.082753863184,30.04719770579885,29.87993040248386,30.56610560735726,30.06053724054328,29.52941175685311,27.93915337040551,27.49183705856683,26.52306883482839,26.21526698621722,25.41927120505517)

for ($n = 0; $n -le 109; $n++) {
    $r = $n + 2
    $lastIdx = 21 - [Math]::Floor($n / 5)
    for ($j = 0; $j -le 21; $j++) {
        $col = $j + 1
        if ($j -le $lastIdx) {
            $ws.Cells.Item($r, $col).Value = $n
        } else {
            $ws.Cells.Item($r, $col).Value = ""
        }
    }
    $rawVal = $rawValues[$n]
    if ($rawVal -eq $null) {
        $ws.Cells.Item($r, 23).Value = ""
    } else {
        $ws.Cells.Item($r, 23).Value = $rawVal
    }
    $finalVal = $finalValues[$n]
    if ($finalVal -eq $null) {
        $ws.Cells.Item($r, 24).Value = ""
    } else {
        $ws.Cells.Item($r, 24).Value = $finalVal
    }
}

# Column A keeps the bold/bordered/centered header-like styling on every
# data row (s="1" in the source); re-apply explicitly for the new rows.
$ws.Range("A2:A111").Font.Bold = $true
$ws.Range("A2:A111").Borders.LineStyle = 1
$ws.Range("A2:A111").HorizontalAlignment = -4108
$ws.Range("A2:A111").VerticalAlignment = -4160
